# Auto-generated Excel COM-interop script
# Updates leve-crafting profit metrics (columns H-N) across multiple worksheets
# to reflect refreshed Market Board price data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 7134.5557
$ws.Cells.Item(116, 9).Value = 1216.6666
$ws.Cells.Item(116, 10).Value = 10093.5
$ws.Cells.Item(116, 11).Value = 1216.6666
$ws.Cells.Item(116, 12).Value = 10093.5
$ws.Cells.Item(116, 13).Value = 2225.3334
$ws.Cells.Item(116, 14).Value = -16977.5
$ws.Cells.Item(135, 8).Value = 944
$ws.Cells.Item(135, 9).Value = 939.8182
$ws.Cells.Item(135, 10).Value = 951.6667
$ws.Cells.Item(135, 11).Value = 8458.363800000001
$ws.Cells.Item(135, 12).Value = 8565.0003
$ws.Cells.Item(135, 13).Value = -5923.363800000001
$ws.Cells.Item(135, 14).Value = -13635.0003
$ws.Cells.Item(137, 8).Value = 2579.9333
$ws.Cells.Item(137, 9).Value = 2269.9
$ws.Cells.Item(137, 10).Value = 3200
$ws.Cells.Item(137, 11).Value = 6809.700000000001
$ws.Cells.Item(137, 12).Value = 9600
$ws.Cells.Item(137, 13).Value = -4259.700000000001
$ws.Cells.Item(137, 14).Value = -14700
$ws.Cells.Item(138, 8).Value = 2073.44
$ws.Cells.Item(138, 9).Value = 1004.5476
$ws.Cells.Item(138, 11).Value = 3013.6428
$ws.Cells.Item(138, 13).Value = 2126.3572

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1782.9259
$ws.Cells.Item(2, 9).Value = 1714.9445
$ws.Cells.Item(2, 10).Value = 1918.8889
$ws.Cells.Item(2, 11).Value = 1714.9445
$ws.Cells.Item(2, 12).Value = 1918.8889
$ws.Cells.Item(2, 13).Value = -1601.9445
$ws.Cells.Item(2, 14).Value = -2144.8889
$ws.Cells.Item(61, 8).Value = 1689.3684
$ws.Cells.Item(61, 9).Value = 1583.2778
$ws.Cells.Item(61, 10).Value = 3599
$ws.Cells.Item(61, 11).Value = 1583.2778
$ws.Cells.Item(61, 12).Value = 3599
$ws.Cells.Item(61, 13).Value = -1371.2778
$ws.Cells.Item(61, 14).Value = -4023
$ws.Cells.Item(74, 8).Value = 93178.37
$ws.Cells.Item(74, 9).Value = 113227.555
$ws.Cells.Item(74, 11).Value = 113227.555
$ws.Cells.Item(74, 13).Value = -112353.555
$ws.Cells.Item(77, 8).Value = 93178.37
$ws.Cells.Item(77, 9).Value = 113227.555
$ws.Cells.Item(77, 11).Value = 566137.7749999999
$ws.Cells.Item(77, 13).Value = -561769.7749999999
$ws.Cells.Item(116, 8).Value = 1782.9259
$ws.Cells.Item(116, 9).Value = 1714.9445
$ws.Cells.Item(116, 10).Value = 1918.8889
$ws.Cells.Item(116, 11).Value = 1714.9445
$ws.Cells.Item(116, 12).Value = 1918.8889
$ws.Cells.Item(116, 13).Value = 579.0554999999999
$ws.Cells.Item(116, 14).Value = -6506.8889
$ws.Cells.Item(132, 8).Value = 1973.6207
$ws.Cells.Item(132, 9).Value = 1834.6111
$ws.Cells.Item(132, 10).Value = 2201.0908
$ws.Cells.Item(132, 11).Value = 5503.8333
$ws.Cells.Item(132, 12).Value = 6603.2724
$ws.Cells.Item(132, 13).Value = -2973.8333
$ws.Cells.Item(132, 14).Value = -11663.2724
$ws.Cells.Item(136, 8).Value = 1689.3684
$ws.Cells.Item(136, 9).Value = 1583.2778
$ws.Cells.Item(136, 10).Value = 3599
$ws.Cells.Item(136, 11).Value = 4749.8334
$ws.Cells.Item(136, 12).Value = 10797
$ws.Cells.Item(136, 13).Value = -2199.8334
$ws.Cells.Item(136, 14).Value = -15897

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1782.9259
$ws.Cells.Item(3, 9).Value = 1714.9445
$ws.Cells.Item(3, 10).Value = 1918.8889
$ws.Cells.Item(3, 11).Value = 1714.9445
$ws.Cells.Item(3, 12).Value = 1918.8889
$ws.Cells.Item(3, 13).Value = -1600.9445
$ws.Cells.Item(3, 14).Value = -2146.8889
$ws.Cells.Item(134, 8).Value = 5945.3716
$ws.Cells.Item(134, 9).Value = 5551.0435
$ws.Cells.Item(134, 10).Value = 6701.1665
$ws.Cells.Item(134, 11).Value = 16653.1305
$ws.Cells.Item(134, 12).Value = 20103.4995
$ws.Cells.Item(134, 13).Value = -14118.1305
$ws.Cells.Item(134, 14).Value = -25173.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 55558024
$ws.Cells.Item(31, 9).Value = 111112680
$ws.Cells.Item(31, 10).Value = 3371.2222
$ws.Cells.Item(31, 11).Value = 111112680
$ws.Cells.Item(31, 12).Value = 3371.2222
$ws.Cells.Item(31, 13).Value = -111112385
$ws.Cells.Item(31, 14).Value = -3961.2222
$ws.Cells.Item(34, 8).Value = 55558024
$ws.Cells.Item(34, 9).Value = 111112680
$ws.Cells.Item(34, 10).Value = 3371.2222
$ws.Cells.Item(34, 11).Value = 111112680
$ws.Cells.Item(34, 12).Value = 3371.2222
$ws.Cells.Item(34, 13).Value = -111112478
$ws.Cells.Item(34, 14).Value = -3775.2222
$ws.Cells.Item(42, 8).Value = 15000
$ws.Cells.Item(42, 10).Value = 15000
$ws.Cells.Item(42, 12).Value = 15000
$ws.Cells.Item(42, 14).Value = -16186
$ws.Cells.Item(105, 8).Value = 1297.375
$ws.Cells.Item(105, 9).Value = 1155.8
$ws.Cells.Item(105, 10).Value = 1533.3334
$ws.Cells.Item(105, 11).Value = 1155.8
$ws.Cells.Item(105, 12).Value = 1533.3334
$ws.Cells.Item(105, 13).Value = 591.2
$ws.Cells.Item(105, 14).Value = -5027.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 715056.5
$ws.Cells.Item(92, 9).Value = 641.2857
$ws.Cells.Item(92, 11).Value = 1923.8571
$ws.Cells.Item(92, 13).Value = -675.8571000000002
$ws.Cells.Item(131, 8).Value = 870.71
$ws.Cells.Item(131, 9).Value = 535.6667
$ws.Cells.Item(131, 10).Value = 916.3977
$ws.Cells.Item(131, 11).Value = 1607.0001
$ws.Cells.Item(131, 12).Value = 2749.1931
$ws.Cells.Item(131, 13).Value = 3432.9999
$ws.Cells.Item(131, 14).Value = -12829.1931

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(33, 8).Value = 7499
$ws.Cells.Item(33, 10).Value = 7499
$ws.Cells.Item(33, 12).Value = 7499
$ws.Cells.Item(33, 14).Value = -8003
$ws.Cells.Item(74, 8).Value = 30000
$ws.Cells.Item(74, 10).Value = 30000
$ws.Cells.Item(74, 12).Value = 30000
$ws.Cells.Item(74, 14).Value = -31872
$ws.Cells.Item(77, 8).Value = 30000
$ws.Cells.Item(77, 10).Value = 30000
$ws.Cells.Item(77, 12).Value = 90000
$ws.Cells.Item(77, 14).Value = -99360
$ws.Cells.Item(97, 8).Value = 1900
$ws.Cells.Item(97, 9).Value = 1750
$ws.Cells.Item(97, 10).Value = 2000
$ws.Cells.Item(97, 11).Value = 1750
$ws.Cells.Item(97, 12).Value = 2000
$ws.Cells.Item(97, 13).Value = -1254
$ws.Cells.Item(97, 14).Value = -2992
$ws.Cells.Item(101, 8).Value = 30000
$ws.Cells.Item(101, 10).Value = 30000
$ws.Cells.Item(101, 12).Value = 30000
$ws.Cells.Item(101, 14).Value = -36490
$ws.Cells.Item(113, 8).Value = 1583.6316
$ws.Cells.Item(113, 9).Value = 1035.1
$ws.Cells.Item(113, 10).Value = 2193.111
$ws.Cells.Item(113, 11).Value = 1035.1
$ws.Cells.Item(113, 12).Value = 2193.111
$ws.Cells.Item(113, 13).Value = 1134.9
$ws.Cells.Item(113, 14).Value = -6533.111

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3610.8333
$ws.Cells.Item(7, 9).Value = 3610.8333
$ws.Cells.Item(7, 11).Value = 3610.8333
$ws.Cells.Item(7, 13).Value = -3498.8333
$ws.Cells.Item(61, 8).Value = 1858.3334
$ws.Cells.Item(61, 9).Value = 1328.5714
$ws.Cells.Item(61, 10).Value = 2600
$ws.Cells.Item(61, 11).Value = 1328.5714
$ws.Cells.Item(61, 12).Value = 2600
$ws.Cells.Item(61, 13).Value = -1126.5714
$ws.Cells.Item(61, 14).Value = -3004
$ws.Cells.Item(93, 8).Value = 70602.39999999999
$ws.Cells.Item(93, 9).Value = 668.8889
$ws.Cells.Item(93, 11).Value = 668.8889
$ws.Cells.Item(93, 13).Value = 579.1111
$ws.Cells.Item(113, 8).Value = 1858.3334
$ws.Cells.Item(113, 9).Value = 1328.5714
$ws.Cells.Item(113, 10).Value = 2600
$ws.Cells.Item(113, 11).Value = 1328.5714
$ws.Cells.Item(113, 12).Value = 2600
$ws.Cells.Item(113, 13).Value = 841.4286
$ws.Cells.Item(113, 14).Value = -6940
$ws.Cells.Item(122, 8).Value = 5379.3
$ws.Cells.Item(122, 9).Value = 6431.75
$ws.Cells.Item(122, 10).Value = 3800.625
$ws.Cells.Item(122, 11).Value = 19295.25
$ws.Cells.Item(122, 12).Value = 11401.875
$ws.Cells.Item(122, 13).Value = -16845.25
$ws.Cells.Item(122, 14).Value = -16301.875
$ws.Cells.Item(126, 8).Value = 3610.8333
$ws.Cells.Item(126, 9).Value = 3610.8333
$ws.Cells.Item(126, 11).Value = 10832.4999
$ws.Cells.Item(126, 13).Value = -8362.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3335755.8
$ws.Cells.Item(132, 9).Value = 4350208
$ws.Cells.Item(132, 10).Value = 2556.6428
$ws.Cells.Item(132, 11).Value = 13050624
$ws.Cells.Item(132, 12).Value = 7669.928400000001
$ws.Cells.Item(132, 13).Value = -13048094
$ws.Cells.Item(132, 14).Value = -12729.9284
